$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix casing of the existing "user" role value -> "User"
$ws.Range("H2").Value = "User"

# Populate the new data row (row 3) with the new ticket/user record
$ws.Range("A3").Value = "nyzopixu@mailinator.com"
$ws.Range("B3").Value = "nyzopixu@mailinator.com"
$ws.Range("C3").Value = "Minim dolore dolore "
$ws.Range("D3").Value = "Adipisicing labore a"
$ws.Range("E3").Value = "VOLUPTATES SUSCIPIT"
$ws.Range("F3").Value = "+1 (165) 371-6338"
$ws.Range("G3").Value = "Excepturi dicta omni"
$ws.Range("H3").Value = "dmine"

# Match the formatting (style + row height) of the row above it
$ws.Range("A2:H2").Copy()
$ws.Range("A3:H3").PasteSpecial(-4122)
$ws.Rows.Item(3).RowHeight = 18.75
